$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Fgf2", "Sdc3", "ECs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 3, 1, 26.12444933333333, 78.373348, 0.7238861157526749, 0.7238861157526749, 16.88264672087378, 151.943820487864, 0.02505985538037673, 0.02505985538037672),
    @("ECs", "Fgf2", "Sdc3", "FAPs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 3, 1, 3.818542, 11.455626, 0.1058085282850919, 0.1058085282850919, 2.467692036385333, 22.209228327468, 0.003662933103887352, 0.003662933103887351),
    @("ECs", "Fgf2", "Sdc3", "sCs", 2, 0.6666666666666666, 0.6462393333333333, 1.938718, 0.03461850536298827, 0.03461850536298827, 3, 1, 6.146179, 18.438537, 0.1703053559622332, 0.1703053559622332, 3.971902619507333, 35.747123575566, 0.005895716878724198, 0.005895716878724196),
    @("FAPs", "Fgf2", "Sdc3", "ECs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 3, 1, 26.12444933333333, 78.373348, 0.7238861157526749, 0.7238861157526749, 400.3003914804826, 3602.703523324344, 0.5941882268261961, 0.594188226826196),
    @("FAPs", "Fgf2", "Sdc3", "FAPs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 3, 1, 3.818542, 11.455626, 0.1058085282850919, 0.1058085282850919, 58.51085463969201, 526.5976917572281, 0.08685092922308321, 0.08685092922308321),
    @("FAPs", "Fgf2", "Sdc3", "sCs", 3, 1, 15.322826, 45.968478, 0.8208310864042159, 0.8208310864042158, 3, 1, 6.146179, 18.438537, 0.1703053559622332, 0.1703053559622332, 94.176831381854, 847.5914824366861, 0.1397919303549366, 0.1397919303549366),
    @("sCs", "Fgf2", "Sdc3", "ECs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 3, 1, 26.12444933333333, 78.373348, 0.7238861157526749, 0.7238861157526749, 70.49390058767466, 634.445105289072, 0.1046380335461021, 0.1046380335461021),
    @("sCs", "Fgf2", "Sdc3", "FAPs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 3, 1, 3.818542, 11.455626, 0.1058085282850919, 0.1058085282850919, 10.303907910296, 92.73517119266401, 0.01529466595812137, 0.01529466595812137),
    @("sCs", "Fgf2", "Sdc3", "sCs", 3, 1, 2.698388, 8.095164, 0.1445504082327959, 0.1445504082327959, 3, 1, 6.146179, 18.438537, 0.1703053559622332, 0.1703053559622332, 16.584775659452, 149.262980935068, 0.02461770872857243, 0.02461770872857243)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $rowVals[$j]
    }
}